$d = $word.ActiveDocument

foreach ($p in $d.Paragraphs) {
    if ($p.Range.Text -like "*Brief Solutions from*anual*") {
        $p.Range.Delete()
        break
    }
}
